$d = $word.ActiveDocument

$d.Content.Find.Execute("13 ± 108", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "13 ± 108 (2)", 2)

$d.Content.Find.Execute("139 ± 258", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "139 ± 258 (74)", 2)

$d.Content.Find.Execute("24 ± 33", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "24 ± 33 (17)", 2)

$d.Content.Find.Execute("93 ± 58", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "93 ± 58 (82)", 2)
